# SystemsAndUsers.xlsx update
# Commit: "Updated the way of fetching the data from the test data files"
#
# The "Type"/"Name" table columns on the Systems_EK1 / Systems_EK2 sheets are
# renamed to "Systemart"/"URL", and the EK2 sheet's URL values are bumped from
# the "...1.toll-collect.de" host names to "...2.toll-collect.de".

$wb = $excel.ActiveWorkbook

# Remember whichever sheet is active right now so we can restore it at the
# end - only the Systems_EK1 / Systems_EK2 cell selections actually change.
$originalActive = $wb.ActiveSheet

# ---- Systems_EK1 (Table1) ----------------------------------------------
$ws1 = $wb.Worksheets.Item("Systems_EK1")

$ws1.Range("A1").Value = "Systemart"
$ws1.Range("B1").Value = "URL"

# Column A/B were manually resized (A keeps its "best fit" auto width, B was
# widened by hand) to better fit the new, longer header text.
$ws1.Columns.Item(1).ColumnWidth = 11.333333333333332
$ws1.Columns.Item(2).ColumnWidth = 49.66666666666667

$ws1.Activate()
$ws1.Range("B23").Select() | Out-Null

# ---- Systems_EK2 (Table15) ---------------------------------------------
$ws2 = $wb.Worksheets.Item("Systems_EK2")

$ws2.Range("A1").Value = "Systemart"
$ws2.Range("B1").Value = "URL"

# Point the EK2 system URLs at the "2" hosts instead of the "1" hosts.
$ws2.Range("B2").Value = "MV2.toll-collect.de"
$ws2.Range("B3").Value = "KO2.toll-collect.de"
$ws2.Range("B4").Value = "IntP2.toll-collect.de"

$ws2.Activate()
$ws2.Range("B2").Select() | Out-Null

# Restore the originally active sheet/tab.
$originalActive.Activate()
